$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels: "<name>_old" -> "<name>_FV2310", "<name>_new" -> "<name>_FV2404" ---
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

$oldCols = @("A","B","C","D","E","F","G","H","I","J")
for ($i = 0; $i -lt $oldCols.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $baseNames[$i] + "_FV2310"
}

# K1 stays "diff" (unchanged)

$newCols = @("L","M","N","O","P","Q","R","S","T","U")
for ($i = 0; $i -lt $newCols.Length; $i++) {
    $ws.Range($newCols[$i] + "1").Value = $baseNames[$i] + "_FV2404"
}

# --- 2. Turn the data range into an Excel Table (ListObject) named Table1 ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U85"), $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (split below row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
